$d = $word.ActiveDocument

# 1) "VIGENCIA DEL PROGRAMA:" -> "VIGENCIA DEL PROGRAMA: 2023"
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("VIGENCIA DEL PROGRAMA:", $true, $false, $false, $false, $false, $true, 1, $false, "VIGENCIA DEL PROGRAMA: 2023", 2)

# 2) "Unidad Nº 1: " -> collapse to single run "Unidad Nº 1: " (fix spellcheck marker split)
for ($i = 1; $i -le 4; $i++) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $r.Find.Execute("Unidad Nº $($i): ", $true, $false, $false, $false, $false, $true, 1, $false, "Unidad Nº $($i): ", 2)
}

# 3) "Envio" -> "Envío"
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("Envio", $true, $false, $false, $false, $false, $true, 1, $false, "Envío", 2)

# 4) Exam paragraph: insert " o entregas" twice
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("examen parcial, teniendo", $true, $false, $false, $false, $false, $true, 1, $false, "examen parcial o entregas, teniendo", 2)

$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("dos parciales en el caso", $true, $false, $false, $false, $false, $true, 1, $false, "dos parciales o entregas en el caso", 2)
